$d = $word.ActiveDocument

# 1. Add a sentence about the unsaved-changes asterisk marker right after the
#    paragraph describing what is stored in the settings folder.
$r = $d.Content
$found = $r.Find.Execute("to load and locations of recently opened files. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $r.Collapse(0)
    $r.InsertAfter("When there are unsaved changes, an asterisk appears next to the file name in the main window. ")
}

# 2. Remove the "unsaved changes...buggy" bullet point from the known-issues list.
$r2 = $d.Content
$found2 = $r2.Find.Execute("The feature checking for unsaved changes is currently a little buggy sometimes. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $para = $r2.Paragraphs(1)
    $para.Range.Delete()
}
